$wb = $excel.ActiveWorkbook

$qa   = $wb.Worksheets.Item("QA")
$prod = $wb.Worksheets.Item("Prod")

# ---------------------------------------------------------------------------
# Prod sheet: column width + row height formatting tweaks
# ---------------------------------------------------------------------------
$prod.Columns.Item(4).ColumnWidth = 31.5

$prod.Rows.Item(33).RowHeight = 90
$prod.Rows.Item(42).RowHeight = 75
$prod.Rows.Item(53).RowHeight = 60

# ---------------------------------------------------------------------------
# QA sheet: update existing cell + append missing input vars
# ---------------------------------------------------------------------------
$qa.Range("C43").Value = "MammothPriceListener,R10PriceService"

$qa.Range("A49").Value = "AdditionalIrmaDeployImplementerRef"
$qa.Range("C49").Value = "."
$qa.Range("C49").WrapText = $true
$qa.Range("F49").Value = "AdditionalIrmaComponentDeploySection"

$qa.Range("A50").Value = "AdditionalIrmaComponentList"
$qa.Range("C50").Value = "*NONE*"
$qa.Range("C50").WrapText = $true

$qa.Range("A51").Value = "SpecialIrmaReleaseDetails"
$qa.Range("C51").Value = "<hr>"
$qa.Range("C51").WrapText = $true

$qa.Range("A52").Value = "SpecialIconReleaseDetails"
$qa.Range("C52").Value = "<hr>"
$qa.Range("C52").WrapText = $true

$qa.Range("A53").Value = "IrmaPostDbUpdateTask"
$qa.Range("A53").Interior.Color = 65535
$qa.Range("C53").Value = "None, continue"
$qa.Range("C53").Interior.Color = 65535
$qa.Range("C53").WrapText = $true

# ---------------------------------------------------------------------------
# Switch the active/selected tab from Prod back to QA, and reset Prod's
# lingering selection back to the top-left cell.
# ---------------------------------------------------------------------------
$prod.Activate()
$prod.Range("A2").Select()
$qa.Activate()
$qa.Range("A2").Select()
